# Update "想去人数" (want-to-go count) values in column F for the
# "展览" and "全部类型" worksheets to match the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 192
    3  = 240
    4  = 261
    5  = 794
    6  = 245
    7  = 5924
    8  = 37
    9  = 67
    10 = 100
    11 = 47
    14 = 180
    15 = 367
    16 = 28
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
